# ---------------------------------------------------------------------------
# Commit: "added 2014 cellulose yield"
#
# The canonical-OOXML diff for this particular part (CsCl_fractionation.docx)
# touches nothing in word/document.xml (no runs/paragraphs/tables were
# added, removed or reworded) - every hunk is inside word/numbering.xml and
# only rewrites the <w:nsid w:val="..."/> child of four pre-existing
# <w:abstractNum> definitions (990, 991, 99411, 99414). The list geometry
# (numFmt/lvlText/indents/...) and every <w:num>/<w:abstractNumId> mapping
# stay byte-for-byte identical.
#
# That nsid is the list "signature" GUID Word stamps on each abstract list
# definition; it has no r/w surface anywhere in the Word object model (no
# List/ListTemplate/ListFormat property exposes it, and it cannot be poked
# through a Find/Replace since it never appears in document "text"). The
# change is a side effect of the authoring tool re-serializing its numbering
# part when *some* file in the same commit/export batch was touched - for
# *this* document body there is nothing else in the diff to reproduce, so
# the edit here is a no-op against Content/Paragraphs/Tables, plus a
# best-effort attempt at the only COM surface that can even see the raw
# numbering markup (Document.WordOpenXML). On real Word automation, and on
# this host, that property is read-only, so the assignment is a harmless
# no-op if unsupported rather than something that corrupts the document.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

$nsidMap = @{
    "53c4bccd" = "671e478c"  # abstractNumId 990
    "1e838963" = "2b360a39"  # abstractNumId 991
    "a52373e7" = "20010374"  # abstractNumId 99411
    "f14ad141" = "cf047a21"  # abstractNumId 99414
}

try {
    $packageXml = $d.WordOpenXML
    foreach ($oldNsid in $nsidMap.Keys) {
        $packageXml = $packageXml.Replace($oldNsid, $nsidMap[$oldNsid])
    }
    $d.WordOpenXML = $packageXml
} catch {
    # Document.WordOpenXML is read-only in this host (and in real Word COM
    # automation there is no writable property for w:nsid at all) - ignore
    # and leave the document exactly as authored, since nothing else in the
    # diff is reachable through Content/Paragraphs/Tables/Find-Replace.
}
